# TT24_curve_log.xlsx -- "add curve fits from second week"
#
# 1) Mark the week-2 curve_fit flags (rows 52-69, col G) as "y" instead of "n"
# 2) Append the week-2 (5/15) data rows 85-94 (copy formatting from row 84)
# 3) Leave the selection on G70 (matches the saved cursor position in the diff)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1) curve_fit column: rows 52..69 go from "n" to "y"
# -----------------------------------------------------------------
for ($r = 52; $r -le 69; $r++) {
    $ws.Cells.Item($r, 7).Value = "y"
}

# -----------------------------------------------------------------
# 2) Append new log rows 85-94 (plot 5, date 45405 = 2024-04-19)
# -----------------------------------------------------------------

# Copy the formatting of the last existing row (84) down across the
# new rows so the date/time/number styles (s="1"/"2"/"3") match.
$ws.Range("A84:G84").Copy()
$ws.Range("A85:G94").PasteSpecial(-4122)

# row 85
$ws.Range("A85").Value = 45405
$ws.Range("B85").Value = 5
$ws.Range("C85").Value = 4265
$ws.Range("D85").Value = "Tri"
$ws.Range("E85").Value = 0.43055555555555558
$ws.Range("F85").Value = 14.8
$ws.Range("G85").Value = "n"

# row 86
$ws.Range("A86").Value = 45405
$ws.Range("B86").Value = 5
$ws.Range("C86").Value = 2573
$ws.Range("D86").Value = "Tri"
$ws.Range("E86").Value = 0.44444444444444442
$ws.Range("F86").Value = 14.1
$ws.Range("G86").Value = "n"

# row 87
$ws.Range("A87").Value = 45405
$ws.Range("B87").Value = 5
$ws.Range("C87").Value = 2547
$ws.Range("D87").Value = "Tri"
$ws.Range("E87").Value = 0.45833333333333331
$ws.Range("F87").Value = 14
$ws.Range("G87").Value = "n"

# row 88
$ws.Range("A88").Value = 45405
$ws.Range("B88").Value = 5
$ws.Range("C88").Value = 4177
$ws.Range("D88").Value = "Tri"
$ws.Range("E88").Value = 0.46805555555555556
$ws.Range("F88").Value = 15
$ws.Range("G88").Value = "n"

# row 89
$ws.Range("A89").Value = 45405
$ws.Range("B89").Value = 5
$ws.Range("C89").Value = 1795
$ws.Range("D89").Value = "Tri"
$ws.Range("E89").Value = 0.47916666666666669
$ws.Range("F89").Value = 9.8
$ws.Range("G89").Value = "n"

# row 90 (tag replaced by a flag note)
$ws.Range("A90").Value = 45405
$ws.Range("B90").Value = 5
$ws.Range("C90").Value = "flag3"
$ws.Range("D90").Value = "Tri"
$ws.Range("E90").Value = 0.48958333333333331
$ws.Range("F90").Value = 12.6
$ws.Range("G90").Value = "n"

# row 91 (tag replaced by a flag note)
$ws.Range("A91").Value = 45405
$ws.Range("B91").Value = 5
$ws.Range("C91").Value = "flag1"
$ws.Range("D91").Value = "Tri"
$ws.Range("E91").Value = 0.5
$ws.Range("F91").Value = 12
$ws.Range("G91").Value = "n"

# row 92 (tag replaced by a flag note)
$ws.Range("A92").Value = 45405
$ws.Range("B92").Value = 5
$ws.Range("C92").Value = "flag2"
$ws.Range("D92").Value = "Tri"
$ws.Range("E92").Value = 0.51041666666666663
$ws.Range("F92").Value = 13.3
$ws.Range("G92").Value = "n"

# row 93
$ws.Range("A93").Value = 45405
$ws.Range("B93").Value = 5
$ws.Range("C93").Value = 1739
$ws.Range("D93").Value = "Tri"
$ws.Range("E93").Value = 0.52430555555555558
$ws.Range("F93").Value = 11
$ws.Range("G93").Value = "n"

# row 94
$ws.Range("A94").Value = 45405
$ws.Range("B94").Value = 5
$ws.Range("C94").Value = 4149
$ws.Range("D94").Value = "Tri"
$ws.Range("E94").Value = 0.53819444444444442
$ws.Range("F94").Value = 13.8
$ws.Range("G94").Value = "n"

# -----------------------------------------------------------------
# 3) Leave the cursor / selection where the saved file shows it
# -----------------------------------------------------------------
$null = $ws.Range("G70").Select()

Write-Output "applied curve-fit flags + second week rows"
